$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '40.178.46'
$ws.Range('E2').Value = '  +0.45%  '
$ws.Range('D3').Value = '2.225.73'
$ws.Range('E3').Value = '  +0.51%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = "'293.84"
$ws.Range('E5').Value = '  +1.58%  '
$ws.Range('D6').Value = "'88.13"
$ws.Range('E6').Value = '  +0.18%  '
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -0.16%  '
$ws.Range('D10').Value = "'30.73"
$ws.Range('E10').Value = '  +0.21%  '
$ws.Range('D11').Value = "'51.01"
$ws.Range('E11').Value = '  +6.74%  '
$ws.Range('D12').Value = "'0.0783"
$ws.Range('E12').Value = '  +0.64%  '
$ws.Range('E13').Value = '  +3.39%  '
$ws.Range('E14').Value = '  -0.19%  '
$ws.Range('E15').Value = '  -1.03%  '
$ws.Range('D16').Value = '2.229.41'
$ws.Range('E16').Value = '  +1.52%  '
$ws.Range('E17').Value = '  +1.18%  '
$ws.Range('D18').Value = '2.030.24'
$ws.Range('E18').Value = '  -20.63%  '
$ws.Range('D19').Value = '40.086.48'
$ws.Range('E19').Value = '  +0.37%  '
$ws.Range('E20').Value = '  +0.60%  '
$ws.Range('D21').Value = "'11.32"
$ws.Range('E21').Value = '  -3.37%  '
$ws.Range('E22').Value = '  -0.35%  '
$ws.Range('D23').Value = "'65.73"
$ws.Range('E23').Value = '  +0.18%  '
$ws.Range('D24').Value = "'236.18"
$ws.Range('E25').Value = '  +0.23%  '
$ws.Range('E26').Value = '  +1.00%  '
$ws.Range('E27').Value = '  -0.17%  '
$ws.Range('D28').Value = "'23.29"
$ws.Range('E28').Value = '  +3.10%  '
$ws.Range('D29').Value = "'9.35"
$ws.Range('E29').Value = '  +1.36%  '
$ws.Range('E30').Value = '  -10.19%  '
$ws.Range('D31').Value = "'159.10"
$ws.Range('E31').Value = '  +4.02%  '
$ws.Range('D32').Value = "'31.98"
$ws.Range('E32').Value = '  -0.22%  '
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('E35').Value = '  +6.51%  '
$ws.Range('D36').Value = "'0.0716"
$ws.Range('E36').Value = '  -0.55%  '
$ws.Range('D37').Value = "'2.33"
$ws.Range('E37').Value = '  -2.59%  '
$ws.Range('E38').Value = '  +1.74%  '
$ws.Range('D40').Value = "'0.0995"
$ws.Range('E40').Value = '  -0.49%  '
$ws.Range('D41').Value = "'15.72"
$ws.Range('E41').Value = '  -0.89%  '
$ws.Range('D42').Value = '2.088.60'
$ws.Range('E42').Value = '  -0.83%  '
$ws.Range('E43').Value = '  -2.62%  '
$ws.Range('D44').Value = "'19.21"
$ws.Range('E44').Value = '  +8.61%  '
$ws.Range('D45').Value = "'10.12"
$ws.Range('E45').Value = '  +1.96%  '
$ws.Range('E46').Value = '  +1.15%  '
$ws.Range('D47').Value = "'2.76"
$ws.Range('E47').Value = '  +2.93%  '
$ws.Range('D48').Value = "'1.90"
$ws.Range('E48').Value = '  -13.32%  '
$ws.Range('D49').Value = '2.438.69'
$ws.Range('E49').Value = '  +0.32%  '
$ws.Range('D50').Value = "'1.47"
$ws.Range('E50').Value = '  +1.36%  '
$ws.Range('E51').Value = '  +3.93%  '
